$d = $word.ActiveDocument

# 1. Merge "...excellent verbal" + " and written communication skills." into one run
$d.Content.Find.Execute(", high level of accuracy and strong attention to detail, with excellent verbal and written communication skills.", $true, $false, $false, $false, $false, $true, 1, $false, ", high level of accuracy and strong attention to detail, with excellent verbal and written communication skills.", 2)

# 2. Change "Python" (Proficient line) to "C,  C++"
$d.Content.Find.Execute("Proficient: Javascript, HTML, MySQL, node.js, PHP, Ruby, Scala, Python", $true, $false, $false, $false, $false, $true, 1, $false, "Proficient: Javascript, HTML, MySQL, node.js, PHP, Ruby, Scala, C,  C++", 2)

# 3. Change "Familiar: C#, C++" to "Familiar: C#, Python"
$d.Content.Find.Execute("Familiar: C#, C++", $true, $false, $false, $false, $false, $true, 1, $false, "Familiar: C#, Python", 2)

# 4. "EMPLOYMENT HISTORY" -> split into "EMPLOYMENT " + "HISTORY"
$d.Content.Find.Execute("EMPLOYMENT HISTORY", $true, $false, $false, $false, $false, $true, 1, $false, "EMPLOYMENT HISTORY", 2)
